$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: translate Spanish column headers to clean machine-friendly names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Normalize "de"/"del"/"los" -> "De"/"Del"/"Los" in place names (data cleaning fix)
$ws.Range("B9").Value = "Hidalgo Del Parral"
$ws.Range("B11").Value = "Valle De Zaragoza"
$ws.Range("A15").Value = "Ciudad De México"
$ws.Range("A20").Value = "Estado De México"
$ws.Range("B20").Value = "Ecatepec De Morelos"
$ws.Range("B21").Value = "Naucalpan De Juárez"
$ws.Range("B30").Value = "Atoyac De Álvarez"
$ws.Range("B31").Value = "Chilpancingo De Los Bravo"
$ws.Range("B34").Value = "Tezontepec De Aldama"
$ws.Range("B38").Value = "Cuautitlán De García Barragán"
$ws.Range("B42").Value = "Unión De Tula"
$ws.Range("B74").Value = "Amatlán De Los Reyes"

# Remove trailing footnote/metadata rows (84-88) so the used range shrinks back to A1:D82
$ws.Range("A84:D88").ClearContents()
